$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_4_9_24"
$ws.Cells.Item(2, 2).Value = 0.9717922210075697
$ws.Cells.Item(2, 3).Value = 0.6443252635344991
$ws.Cells.Item(2, 4).Value = 0.9903848125876391
$ws.Cells.Item(2, 5).Value = 0.1263598138728296
$ws.Cells.Item(2, 6).Value = 0.8565089306939007
$ws.Cells.Item(2, 7).Value = 0.1886255480072607
$ws.Cells.Item(2, 8).Value = 2.378398600476376
$ws.Cells.Item(2, 9).Value = 0.03479649045852942
$ws.Cells.Item(2, 10).Value = 1.143351198438473
$ws.Cells.Item(2, 11).Value = 0.589074367756487
$ws.Cells.Item(2, 12).Value = 0.9367030350697056
$ws.Cells.Item(2, 13).Value = 0.4343104281585473
$ws.Cells.Item(2, 14).Value = 1.01829693772482
$ws.Cells.Item(2, 15).Value = 0.4527999222317838
$ws.Cells.Item(2, 16).Value = 125.3359829132492
$ws.Cells.Item(2, 17).Value = 199.6874082302094

$ws.Cells.Item(3, 1).Value = "model_4_9_23"
$ws.Cells.Item(3, 2).Value = 0.9718052954100157
$ws.Cells.Item(3, 3).Value = 0.6442025115762624
$ws.Cells.Item(3, 4).Value = 0.9904431072102927
$ws.Cells.Item(3, 5).Value = 0.1267520321743291
$ws.Cells.Item(3, 6).Value = 0.8565971706638824
$ws.Cells.Item(3, 7).Value = 0.1885381194179025
$ws.Cells.Item(3, 8).Value = 2.379219443386326
$ws.Cells.Item(3, 9).Value = 0.03458552751064734
$ws.Cells.Item(3, 10).Value = 1.142837894137469
$ws.Cells.Item(3, 11).Value = 0.5887121159119699
$ws.Cells.Item(3, 12).Value = 0.9375285186510997
$ws.Cells.Item(3, 13).Value = 0.4342097643051139
$ws.Cells.Item(3, 14).Value = 1.018288457031341
$ws.Cells.Item(3, 15).Value = 0.4526949729097068
$ws.Cells.Item(3, 16).Value = 125.3369101350203
$ws.Cells.Item(3, 17).Value = 199.6883354519805

$ws.Cells.Item(4, 1).Value = "model_4_9_22"
$ws.Cells.Item(4, 2).Value = 0.9718187554754496
$ws.Cells.Item(4, 3).Value = 0.6440651465571539
$ws.Cells.Item(4, 4).Value = 0.9905075492532195
$ws.Cells.Item(4, 5).Value = 0.1271721920224395
$ws.Cells.Item(4, 6).Value = 0.8566926540693462
$ws.Cells.Item(4, 7).Value = 0.1884481118983677
$ws.Cells.Item(4, 8).Value = 2.380138003901612
$ws.Cells.Item(4, 9).Value = 0.03435231760680846
$ws.Cells.Item(4, 10).Value = 1.142288022149548
$ws.Cells.Item(4, 11).Value = 0.588320127567491
$ws.Cells.Item(4, 12).Value = 0.9384430550965113
$ws.Cells.Item(4, 13).Value = 0.4341061067277995
$ws.Cells.Item(4, 14).Value = 1.018279726178087
$ws.Cells.Item(4, 15).Value = 0.4525869024147254
$ws.Cells.Item(4, 16).Value = 125.3378651568621
$ws.Cells.Item(4, 17).Value = 199.6892904738224

$ws.Cells.Item(5, 1).Value = "model_4_9_21"
$ws.Cells.Item(5, 2).Value = 0.9718323715835391
$ws.Cells.Item(5, 3).Value = 0.6439111966632446
$ws.Cells.Item(5, 4).Value = 0.9905789662238779
$ws.Cells.Item(5, 5).Value = 0.1276176033321479
$ws.Cells.Item(5, 6).Value = 0.856795064649758
$ws.Cells.Item(5, 7).Value = 0.1883570609208774
$ws.Cells.Item(5, 8).Value = 2.381167467551058
$ws.Cells.Item(5, 9).Value = 0.03409386607263494
$ws.Cells.Item(5, 10).Value = 1.141705103045277
$ws.Cells.Item(5, 11).Value = 0.5878997010684804
$ws.Cells.Item(5, 12).Value = 0.9394573492878925
$ws.Cells.Item(5, 13).Value = 0.4340012222573543
$ws.Cells.Item(5, 14).Value = 1.018270894107975
$ws.Cells.Item(5, 15).Value = 0.4524775527952325
$ws.Cells.Item(5, 16).Value = 125.3388317144981
$ws.Cells.Item(5, 17).Value = 199.6902570314583

$ws.Cells.Item(6, 1).Value = "model_4_9_20"
$ws.Cells.Item(6, 2).Value = 0.9718457999013318
$ws.Cells.Item(6, 3).Value = 0.6437383655564448
$ws.Cells.Item(6, 4).Value = 0.9906580030985359
$ws.Cells.Item(6, 5).Value = 0.1280853617145294
$ws.Cells.Item(6, 6).Value = 0.8569044297485648
$ws.Cells.Item(6, 7).Value = 0.1882672656979653
$ws.Cells.Item(6, 8).Value = 2.382323189958042
$ws.Cells.Item(6, 9).Value = 0.03380783879755799
$ws.Cells.Item(6, 10).Value = 1.141092937859233
$ws.Cells.Item(6, 11).Value = 0.5874507241617932
$ws.Cells.Item(6, 12).Value = 0.9405845480489745
$ws.Cells.Item(6, 13).Value = 0.4338977594986695
$ws.Cells.Item(6, 14).Value = 1.018262183847785
$ws.Cells.Item(6, 15).Value = 0.4523696854127131
$ws.Cells.Item(6, 16).Value = 125.3397853993056
$ws.Cells.Item(6, 17).Value = 199.6912107162659

$ws.Cells.Item(7, 1).Value = "model_4_9_19"
$ws.Cells.Item(7, 2).Value = 0.9718587903295411
$ws.Cells.Item(7, 3).Value = 0.6435442626850616
$ws.Cells.Item(7, 4).Value = 0.9907453497999491
$ws.Cells.Item(7, 5).Value = 0.128574895892186
$ws.Cells.Item(7, 6).Value = 0.8570211582143188
$ws.Cells.Item(7, 7).Value = 0.1881803986447141
$ws.Cells.Item(7, 8).Value = 2.383621156752743
$ws.Cells.Item(7, 9).Value = 0.03349173901375161
$ws.Cells.Item(7, 10).Value = 1.140452274234107
$ws.Cells.Item(7, 11).Value = 0.586971518400099
$ws.Cells.Item(7, 12).Value = 0.9418347843870404
$ws.Cells.Item(7, 13).Value = 0.4337976471175404
$ws.Cells.Item(7, 14).Value = 1.018253757624081
$ws.Cells.Item(7, 15).Value = 0.4522653110402581
$ws.Cells.Item(7, 16).Value = 125.3407084179669
$ws.Cells.Item(7, 17).Value = 199.6921337349271

$ws.Cells.Item(8, 1).Value = "model_4_9_18"
$ws.Cells.Item(8, 2).Value = 0.9718706674181066
$ws.Cells.Item(8, 3).Value = 0.643326114734452
$ws.Cells.Item(8, 4).Value = 0.9908418172217829
$ws.Cells.Item(8, 5).Value = 0.1290793538811996
$ws.Cells.Item(8, 6).Value = 0.8571439354637506
$ws.Cells.Item(8, 7).Value = 0.1881009764987884
$ws.Cells.Item(8, 8).Value = 2.385079913102948
$ws.Cells.Item(8, 9).Value = 0.03314263217064601
$ws.Cells.Item(8, 10).Value = 1.139792079504675
$ws.Cells.Item(8, 11).Value = 0.5864674805464984
$ws.Cells.Item(8, 12).Value = 0.9432282511173947
$ws.Cells.Item(8, 13).Value = 0.4337060946064609
$ws.Cells.Item(8, 14).Value = 1.018246053566634
$ws.Cells.Item(8, 15).Value = 0.4521698609492421
$ws.Cells.Item(8, 16).Value = 125.341552702615
$ws.Cells.Item(8, 17).Value = 199.6929780195753

$ws.Cells.Item(9, 1).Value = "model_4_9_17"
$ws.Cells.Item(9, 2).Value = 0.9718807421162085
$ws.Cells.Item(9, 3).Value = 0.6430806151365382
$ws.Cells.Item(9, 4).Value = 0.9909483139474888
$ws.Cells.Item(9, 5).Value = 0.1295891574037478
$ws.Cells.Item(9, 6).Value = 0.857272225362482
$ws.Cells.Item(9, 7).Value = 0.1880336069461911
$ws.Cells.Item(9, 8).Value = 2.386721570044619
$ws.Cells.Item(9, 9).Value = 0.03275723018720361
$ws.Cells.Item(9, 10).Value = 1.139124888963616
$ws.Cells.Item(9, 11).Value = 0.5859408115952512
$ws.Cells.Item(9, 12).Value = 0.9447743243316306
$ws.Cells.Item(9, 13).Value = 0.4336284203626316
$ws.Cells.Item(9, 14).Value = 1.018239518627324
$ws.Cells.Item(9, 15).Value = 0.4520888799520452
$ws.Cells.Item(9, 16).Value = 125.3422691435517
$ws.Cells.Item(9, 17).Value = 199.6936944605119

$ws.Cells.Item(10, 1).Value = "model_4_9_16"
$ws.Cells.Item(10, 2).Value = 0.971888124201553
$ws.Cells.Item(10, 3).Value = 0.6428038823709414
$ws.Cells.Item(10, 4).Value = 0.9910658176325603
$ws.Cells.Item(10, 5).Value = 0.1300975116622255
$ws.Cells.Item(10, 6).Value = 0.857404879650878
$ws.Cells.Item(10, 7).Value = 0.1879842429075006
$ws.Cells.Item(10, 8).Value = 2.388572083322402
$ws.Cells.Item(10, 9).Value = 0.03233199501693813
$ws.Cells.Item(10, 10).Value = 1.138459595104781
$ws.Cells.Item(10, 11).Value = 0.5853962255004862
$ws.Cells.Item(10, 12).Value = 0.9464863476185442
$ws.Cells.Item(10, 13).Value = 0.4335714968808496
$ws.Cells.Item(10, 14).Value = 1.018234730247641
$ws.Cells.Item(10, 15).Value = 0.4520295331197958
$ws.Cells.Item(10, 16).Value = 125.3427942679732
$ws.Cells.Item(10, 17).Value = 199.6942195849335

$ws.Cells.Item(11, 1).Value = "model_4_9_15"
$ws.Cells.Item(11, 2).Value = 0.9718915472140341
$ws.Cells.Item(11, 3).Value = 0.6424916103141829
$ws.Cells.Item(11, 4).Value = 0.991195317124005
$ws.Cells.Item(11, 5).Value = 0.1305851127709732
$ws.Cells.Item(11, 6).Value = 0.8575397157468688
$ws.Cells.Item(11, 7).Value = 0.1879613532072782
$ws.Cells.Item(11, 8).Value = 2.390660248003828
$ws.Cells.Item(11, 9).Value = 0.03186334811228722
$ws.Cells.Item(11, 10).Value = 1.137821461327398
$ws.Cells.Item(11, 11).Value = 0.5848426824236895
$ws.Cells.Item(11, 12).Value = 0.9483924569323002
$ws.Cells.Item(11, 13).Value = 0.4335450993925294
$ws.Cells.Item(11, 14).Value = 1.018232509915221
$ws.Cells.Item(11, 15).Value = 0.4520020118357476
$ws.Cells.Item(11, 16).Value = 125.3430378106612
$ws.Cells.Item(11, 17).Value = 199.6944631276214

$ws.Cells.Item(12, 1).Value = "model_4_9_14"
$ws.Cells.Item(12, 2).Value = 0.9718894944306053
$ws.Cells.Item(12, 3).Value = 0.6421387557653326
$ws.Cells.Item(12, 4).Value = 0.9913380836504287
$ws.Cells.Item(12, 5).Value = 0.1310342982014479
$ws.Cells.Item(12, 6).Value = 0.8576742770880805
$ws.Cells.Item(12, 7).Value = 0.1879750801795181
$ws.Cells.Item(12, 8).Value = 2.393019787996735
$ws.Cells.Item(12, 9).Value = 0.03134668901231864
$ws.Cells.Item(12, 10).Value = 1.13723360295228
$ws.Cells.Item(12, 11).Value = 0.5842902672986082
$ws.Cells.Item(12, 12).Value = 0.9505111534356288
$ws.Cells.Item(12, 13).Value = 0.4335609301811201
$ws.Cells.Item(12, 14).Value = 1.018233841450418
$ws.Cells.Item(12, 15).Value = 0.4520185165737828
$ws.Cells.Item(12, 16).Value = 125.3428917543492
$ws.Cells.Item(12, 17).Value = 199.6943170713095

$ws.Cells.Item(13, 1).Value = "model_4_9_13"
$ws.Cells.Item(13, 2).Value = 0.9718800262165423
$ws.Cells.Item(13, 3).Value = 0.6417393812581498
$ws.Cells.Item(13, 4).Value = 0.9914950937930742
$ws.Cells.Item(13, 5).Value = 0.131420574427667
$ws.Cells.Item(13, 6).Value = 0.8578050805607454
$ws.Cells.Item(13, 7).Value = 0.1880383941705546
$ws.Cells.Item(13, 8).Value = 2.395690407165217
$ws.Cells.Item(13, 9).Value = 0.03077848355816053
$ws.Cells.Item(13, 10).Value = 1.136728075169574
$ws.Cells.Item(13, 11).Value = 0.5837532793638674
$ws.Cells.Item(13, 12).Value = 0.9528504050068313
$ws.Cells.Item(13, 13).Value = 0.4336339402889891
$ws.Cells.Item(13, 14).Value = 1.018239982994675
$ws.Cells.Item(13, 15).Value = 0.4520946348730955
$ws.Cells.Item(13, 16).Value = 125.3422182253833
$ws.Cells.Item(13, 17).Value = 199.6936435423436

$ws.Cells.Item(14, 1).Value = "model_4_9_12"
$ws.Cells.Item(14, 2).Value = 0.9718605469321667
$ws.Cells.Item(14, 3).Value = 0.6412865853495553
$ws.Cells.Item(14, 4).Value = 0.9916674641454583
$ws.Cells.Item(14, 5).Value = 0.1317079538919895
$ws.Cells.Item(14, 6).Value = 0.8579269458336721
$ws.Cells.Item(14, 7).Value = 0.1881686522348682
$ws.Cells.Item(14, 8).Value = 2.398718255490922
$ws.Cells.Item(14, 9).Value = 0.03015469089923095
$ws.Cells.Item(14, 10).Value = 1.136351975649248
$ws.Cells.Item(14, 11).Value = 0.5832529854504701
$ws.Cells.Item(14, 12).Value = 0.9554567279476791
$ws.Cells.Item(14, 13).Value = 0.4337841078634258
$ws.Cells.Item(14, 14).Value = 1.018252618206162
$ws.Cells.Item(14, 15).Value = 0.4522511953920657
$ws.Cells.Item(14, 16).Value = 125.340833263876
$ws.Cells.Item(14, 17).Value = 199.6922585808362

$ws.Cells.Item(15, 1).Value = "model_4_9_11"
$ws.Cells.Item(15, 2).Value = 0.9718278498294214
$ws.Cells.Item(15, 3).Value = 0.6407723519470817
$ws.Cells.Item(15, 4).Value = 0.9918567942627114
$ws.Cells.Item(15, 5).Value = 0.1318519661774877
$ws.Cells.Item(15, 6).Value = 0.858033378331996
$ws.Cells.Item(15, 7).Value = 0.1883872979114829
$ws.Cells.Item(15, 8).Value = 2.40215693662109
$ws.Cells.Item(15, 9).Value = 0.02946952239070585
$ws.Cells.Item(15, 10).Value = 1.136163503756782
$ws.Cells.Item(15, 11).Value = 0.5828160477583749
$ws.Cells.Item(15, 12).Value = 0.9583470991455815
$ws.Cells.Item(15, 13).Value = 0.4340360560039718
$ws.Cells.Item(15, 14).Value = 1.018273827137673
$ws.Cells.Item(15, 15).Value = 0.4525138694865593
$ws.Cells.Item(15, 16).Value = 125.3385106799295
$ws.Cells.Item(15, 17).Value = 199.6899359968898

$ws.Cells.Item(16, 1).Value = "model_4_9_10"
$ws.Cells.Item(16, 2).Value = 0.9717779005511916
$ws.Cells.Item(16, 3).Value = 0.6401871006266241
$ws.Cells.Item(16, 4).Value = 0.9920643179748591
$ws.Cells.Item(16, 5).Value = 0.1317964460854226
$ws.Cells.Item(16, 6).Value = 0.8581159402075971
$ws.Cells.Item(16, 7).Value = 0.1887213089650008
$ws.Cells.Item(16, 8).Value = 2.406070514895824
$ws.Cells.Item(16, 9).Value = 0.02871851291372167
$ws.Cells.Item(16, 10).Value = 1.136236164063402
$ws.Cells.Item(16, 11).Value = 0.5824771062137498
$ws.Cells.Item(16, 12).Value = 0.9615543818681257
$ws.Cells.Item(16, 13).Value = 0.4344206589988565
$ws.Cells.Item(16, 14).Value = 1.018306226669497
$ws.Cells.Item(16, 15).Value = 0.4529148458271742
$ws.Cells.Item(16, 16).Value = 125.3349678157809
$ws.Cells.Item(16, 17).Value = 199.6863931327411

$ws.Cells.Item(17, 1).Value = "model_4_9_9"
$ws.Cells.Item(17, 2).Value = 0.9717055995205491
$ws.Cells.Item(17, 3).Value = 0.6395198508497292
$ws.Cells.Item(17, 4).Value = 0.9922911149307033
$ws.Cells.Item(17, 5).Value = 0.1314700494671966
$ws.Cells.Item(17, 6).Value = 0.858163773645835
$ws.Cells.Item(17, 7).Value = 0.1892047862898225
$ws.Cells.Item(17, 8).Value = 2.410532417226324
$ws.Cells.Item(17, 9).Value = 0.02789775531726428
$ws.Cells.Item(17, 10).Value = 1.136663326149741
$ws.Cells.Item(17, 11).Value = 0.5822807354394304
$ws.Cells.Item(17, 12).Value = 0.9651162927913159
$ws.Cells.Item(17, 13).Value = 0.4349767652298482
$ws.Cells.Item(17, 14).Value = 1.018353124635319
$ws.Cells.Item(17, 15).Value = 0.4534946266517178
$ws.Cells.Item(17, 16).Value = 125.3298506503388
$ws.Cells.Item(17, 17).Value = 199.6812759672991

$ws.Cells.Item(18, 1).Value = "model_4_9_8"
$ws.Cells.Item(18, 2).Value = 0.971604463148229
$ws.Cells.Item(18, 3).Value = 0.6387572958164116
$ws.Cells.Item(18, 4).Value = 0.9925387973594993
$ws.Cells.Item(18, 5).Value = 0.1307773442693168
$ws.Cells.Item(18, 6).Value = 0.858162508098249
$ws.Cells.Item(18, 7).Value = 0.1898810856772172
$ws.Cells.Item(18, 8).Value = 2.415631626245361
$ws.Cells.Item(18, 9).Value = 0.02700141509000379
$ws.Cells.Item(18, 10).Value = 1.137569883941767
$ws.Cells.Item(18, 11).Value = 0.5822859308961763
$ws.Cells.Item(18, 12).Value = 0.969058922090095
$ws.Cells.Item(18, 13).Value = 0.4357534689216108
$ws.Cells.Item(18, 14).Value = 1.018418726606554
$ws.Cells.Item(18, 15).Value = 0.454304396227637
$ws.Cells.Item(18, 16).Value = 125.3227145352271
$ws.Cells.Item(18, 17).Value = 199.6741398521874

$ws.Cells.Item(19, 1).Value = "model_4_9_7"
$ws.Cells.Item(19, 2).Value = 0.9714663956785866
$ws.Cells.Item(19, 3).Value = 0.6378840288169932
$ws.Cells.Item(19, 4).Value = 0.992808492958439
$ws.Cells.Item(19, 5).Value = 0.1296008345731494
$ws.Cells.Item(19, 6).Value = 0.8580939192309756
$ws.Cells.Item(19, 7).Value = 0.1908043434824586
$ws.Cells.Item(19, 8).Value = 2.421471166691497
$ws.Cells.Item(19, 9).Value = 0.02602541119816641
$ws.Cells.Item(19, 10).Value = 1.139109606807596
$ws.Cells.Item(19, 11).Value = 0.5825675090028813
$ws.Cells.Item(19, 12).Value = 0.9734369958266043
$ws.Cells.Item(19, 13).Value = 0.4368115651885359
$ws.Cells.Item(19, 14).Value = 1.01850828388416
$ws.Cells.Item(19, 15).Value = 0.4554075378432065
$ws.Cells.Item(19, 16).Value = 125.3130135112879
$ws.Cells.Item(19, 17).Value = 199.6644388282481

$ws.Cells.Item(20, 1).Value = "model_4_9_6"
$ws.Cells.Item(20, 2).Value = 0.9712813709844371
$ws.Cells.Item(20, 3).Value = 0.6368815474838697
$ws.Cells.Item(20, 4).Value = 0.9931013379468258
$ws.Cells.Item(20, 5).Value = 0.1277942052115866
$ws.Cells.Item(20, 6).Value = 0.85793500153603
$ws.Cells.Item(20, 7).Value = 0.1920416044641965
$ws.Cells.Item(20, 8).Value = 2.428174763982096
$ws.Cells.Item(20, 9).Value = 0.02496563176722902
$ws.Cells.Item(20, 10).Value = 1.141473980469061
$ws.Cells.Item(20, 11).Value = 0.5832199143485803
$ws.Cells.Item(20, 12).Value = 0.9782820155402463
$ws.Cells.Item(20, 13).Value = 0.4382255178149676
$ws.Cells.Item(20, 14).Value = 1.018628299901987
$ws.Cells.Item(20, 15).Value = 0.4568816853602308
$ws.Cells.Item(20, 16).Value = 125.3000864810211
$ws.Cells.Item(20, 17).Value = 199.6515117979813

$ws.Cells.Item(21, 1).Value = "model_4_9_5"
$ws.Cells.Item(21, 2).Value = 0.97103677273585
$ws.Cells.Item(21, 3).Value = 0.6357278327555704
$ws.Cells.Item(21, 4).Value = 0.9934185766658529
$ws.Cells.Item(21, 5).Value = 0.1251703143129216
$ws.Cells.Item(21, 6).Value = 0.8576566475618153
$ws.Cells.Item(21, 7).Value = 0.19367723407877
$ws.Cells.Item(21, 8).Value = 2.435889659682601
$ws.Cells.Item(21, 9).Value = 0.02381757363936419
$ws.Cells.Item(21, 10).Value = 1.14490792141088
$ws.Cells.Item(21, 11).Value = 0.5843626418518738
$ws.Cells.Item(21, 12).Value = 0.9836579343485444
$ws.Cells.Item(21, 13).Value = 0.4400877572470859
$ws.Cells.Item(21, 14).Value = 1.018786958225395
$ws.Cells.Item(21, 15).Value = 0.4588232041803414
$ws.Cells.Item(21, 16).Value = 125.2831244946727
$ws.Cells.Item(21, 17).Value = 199.6345498116329

$ws.Cells.Item(22, 1).Value = "model_4_9_4"
$ws.Cells.Item(22, 2).Value = 0.9707167015997912
$ws.Cells.Item(22, 3).Value = 0.6343965840285308
$ws.Cells.Item(22, 4).Value = 0.9937608917296096
$ws.Cells.Item(22, 5).Value = 0.121491154213103
$ws.Cells.Item(22, 6).Value = 0.8572211412657819
$ws.Cells.Item(22, 7).Value = 0.1958175512393869
$ws.Cells.Item(22, 8).Value = 2.444791726050151
$ws.Cells.Item(22, 9).Value = 0.02257876649614553
$ws.Cells.Item(22, 10).Value = 1.149722915244924
$ws.Cells.Item(22, 11).Value = 0.5861505273086512
$ws.Cells.Item(22, 12).Value = 0.9896114442189643
$ws.Cells.Item(22, 13).Value = 0.4425127695777681
$ws.Cells.Item(22, 14).Value = 1.018994571935271
$ws.Cells.Item(22, 15).Value = 0.4613514543064083
$ws.Cells.Item(22, 16).Value = 125.2611438285907
$ws.Cells.Item(22, 17).Value = 199.6125691455509

$ws.Cells.Item(23, 1).Value = "model_4_9_3"
$ws.Cells.Item(23, 2).Value = 0.9703017710416699
$ws.Cells.Item(23, 3).Value = 0.632856504787862
$ws.Cells.Item(23, 4).Value = 0.9941287509646656
$ws.Cells.Item(23, 5).Value = 0.1164698129808049
$ws.Cells.Item(23, 6).Value = 0.8565827955838209
$ws.Cells.Item(23, 7).Value = 0.1985921937921232
$ws.Cells.Item(23, 8).Value = 2.455090243023917
$ws.Cells.Item(23, 9).Value = 0.02124751731568201
$ws.Cells.Item(23, 10).Value = 1.156294449621297
$ws.Cells.Item(23, 11).Value = 0.5887711299763271
$ws.Cells.Item(23, 12).Value = 0.9962031128969194
$ws.Cells.Item(23, 13).Value = 0.4456368407034176
$ws.Cells.Item(23, 14).Value = 1.019263716081079
$ws.Cells.Item(23, 15).Value = 0.4646085236075955
$ws.Cells.Item(23, 16).Value = 125.2330036686525
$ws.Cells.Item(23, 17).Value = 199.5844289856127

$ws.Cells.Item(24, 1).Value = "model_4_9_2"
$ws.Cells.Item(24, 2).Value = 0.9697676525531529
$ws.Cells.Item(24, 3).Value = 0.6310701258151239
$ws.Cells.Item(24, 4).Value = 0.9945219994535051
$ws.Cells.Item(24, 5).Value = 0.1097438464094224
$ws.Cells.Item(24, 6).Value = 0.855683979779948
$ws.Cells.Item(24, 7).Value = 0.2021638465842265
$ws.Cells.Item(24, 8).Value = 2.467035767440135
$ws.Cells.Item(24, 9).Value = 0.01982438673040152
$ws.Cells.Item(24, 10).Value = 1.165096862859793
$ws.Cells.Item(24, 11).Value = 0.5924610415084967
$ws.Cells.Item(24, 12).Value = 1.003499272764066
$ws.Cells.Item(24, 13).Value = 0.4496263410702563
$ws.Cells.Item(24, 14).Value = 1.019610171316874
$ws.Cells.Item(24, 15).Value = 0.4687678652644547
$ws.Cells.Item(24, 16).Value = 125.1973535773165
$ws.Cells.Item(24, 17).Value = 199.5487788942768

$ws.Cells.Item(25, 1).Value = "model_4_9_1"
$ws.Cells.Item(25, 2).Value = 0.9690842701259327
$ws.Cells.Item(25, 3).Value = 0.6289922275627033
$ws.Cells.Item(25, 4).Value = 0.994940187524691
$ws.Cells.Item(25, 5).Value = 0.100861673024114
$ws.Cells.Item(25, 6).Value = 0.8544526598589313
$ws.Cells.Item(25, 7).Value = 0.206733628021735
$ws.Cells.Item(25, 8).Value = 2.480930682621914
$ws.Cells.Item(25, 9).Value = 0.01831100206041728
$ws.Cells.Item(25, 10).Value = 1.176721149088944
$ws.Cells.Item(25, 11).Value = 0.5975159833072203
$ws.Cells.Item(25, 12).Value = 1.011571102265411
$ws.Cells.Item(25, 13).Value = 0.4546796982731195
$ws.Cells.Item(25, 14).Value = 1.0200534464048
$ws.Cells.Item(25, 15).Value = 0.4740363543453354
$ws.Cells.Item(25, 16).Value = 125.1526482710662
$ws.Cells.Item(25, 17).Value = 199.5040735880265

$ws.Cells.Item(26, 1).Value = "model_4_9_0"
$ws.Cells.Item(26, 2).Value = 0.9682146042791846
$ws.Cells.Item(26, 3).Value = 0.6265687792170147
$ws.Cells.Item(26, 4).Value = 0.9953818892039651
$ws.Cells.Item(26, 5).Value = 0.08927712869585192
$ws.Cells.Item(26, 6).Value = 0.8528008158004662
$ws.Cells.Item(26, 7).Value = 0.2125490875433826
$ws.Cells.Item(26, 8).Value = 2.497136292868486
$ws.Cells.Item(26, 9).Value = 0.01671252377713197
$ws.Cells.Item(26, 10).Value = 1.191882084736602
$ws.Cells.Item(26, 11).Value = 0.6042973042568668
$ws.Cells.Item(26, 12).Value = 1.020486367136452
$ws.Cells.Item(26, 13).Value = 0.4610304627065142
$ws.Cells.Item(26, 14).Value = 1.020617553981069
$ws.Cells.Item(26, 15).Value = 0.4806574839685543
$ws.Cells.Item(26, 16).Value = 125.0971646341799
$ws.Cells.Item(26, 17).Value = 199.4485899511401
